$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 6; this pushes the existing rows 6-13 down to 7-14,
# carrying their values/formatting with them.
$ws.Rows.Item(6).Insert()

# Fill in the new row 6 with the "water" resource entry (raw material model).
$ws.Cells.Item(6, 1).Value = "water"
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = "rawMaterialModel"
$ws.Cells.Item(6, 4).Value = $null
$ws.Cells.Item(6, 5).Value = 10
$ws.Cells.Item(6, 6).Value = 100

# Match the formatting used by the other raw-material rows (rows 2-5).
$ws.Range("A6:F6").Style = $ws.Range("A5:F5").Style

# Update the selection to match the saved selection state.
$ws.Range("H21").Select()
